$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (and update the matching "Instance" cell value used in H3)
$ws.Name = "Sedan_HambaLG"
$ws.Range("H3").Value = "Sedan_HambaLG"

# Update the frozen-pane scroll position and active selection
$ws.Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("D3").Select()
$win.FreezePanes = $true
$ws.Range("H4").Select()
